$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7:AF12").ClearContents()

$ws.Rows("7:12").Select() | Out-Null
